$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.175.95"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").Value = "1.901.45"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9992"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5202"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3765"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07266"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.07%  "
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9037"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08324"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.88%  "
$ws.Range("D13").Value = "1.923.82"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "96.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.292"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008650"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.17%  "
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9991"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "27.213.74"
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.085"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("D22").Value = "2.155.91"
$ws.Range("E22").Value = "  +1.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.443"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.325"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Value = "  +1.33%  "
$ws.Range("E28").Value = "  +1.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.826"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09252"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05072"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7979"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.246"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.416"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.946"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.598"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5732"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02003"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.078"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.030"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.603"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "116.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.32%  "
$ws.Range("E45").Value = "  +1.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4861"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9990"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.631"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.39%  "
